$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11 entered first (test case 3b, short session)
$ws.Range("A11").Value = "Wednesday (2/27/19) 6:00 PM - 3:30 AM"
$ws.Range("B11").Value = "Fixing sll instruction, get Project 1 hamming weight and distance to work"
$ws.Range("C11").Value = "Continue the activity"
$ws.Range("D11").Value = "Richard, Syed"

# New row 10 entered second (test case 3b continuation)
$ws.Range("A10").Value = "Wednesday (2/27/19) 6:00 PM - 10:30 PM"
$ws.Range("B10").Value = "Fixing sll instruction, get Project 1 hamming weight and distance to work"
$ws.Range("C10").Value = "Continue the activity"
$ws.Range("D10").Value = "Richard, Syed, Rami"

# New row 12 entered third (test case 4b, project turn-in)
$ws.Range("A12").Value = "Thursday (2/28/19) 6:00 PM - 11:59 PM"
$ws.Range("B12").Value = "Finish the project, clean source folders, write lab reports"
$ws.Range("C12").Value = "Turn in project"
$ws.Range("D12").Value = "Richard, Syed"

$ws.Range("D12").Select()
